# Actualización automática del tracker
# Agrega la fila 54 con el nuevo resultado pendiente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$row = $lastRow + 1

$ws.Cells.Item($row, 1).Value = 14581503

# La columna "fecha" se guarda como texto (p.ej. "2025-09-01"), no como
# fecha real, así que forzamos formato de texto antes de escribirla para
# que Excel no la auto-convierta a un número de serie de fecha.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2025-09-01"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "Sumit Nagal"
$ws.Cells.Item($row, 4).Value = "Mili Poljičak"
$ws.Cells.Item($row, 5).Value = "Gana Mili Poljičak"
$ws.Cells.Item($row, 6).Value = 2.25

# "resultado" y "profit" quedan vacíos porque el partido todavía no se jugó;
# se tocan con .Style para que la celda exista (vacía) igual que en el resto
# de filas pendientes del tracker.
$ws.Cells.Item($row, 7).Style = "Normal"
$ws.Cells.Item($row, 8).Style = "Normal"
